# Apply updated crypto price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.715.46'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '1.888.48'
$ws.Range("E3").Value = '  +1.01%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("D5").Value = '''248.10'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '''0.4737'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.2920'
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").Value = '''0.06530'
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("D10").Value = '''22.00'
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("D11").Value = '''0.07804'
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.891.48'
$ws.Range("E12").Value = '  +1.03%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '''96.76'
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").Value = '''0.7360'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").Value = '''5.248'
$ws.Range("E15").Value = '  +2.61%  '
$ws.Range("D16").Value = '''283.57'
$ws.Range("E16").Value = '  +3.72%  '
$ws.Range("D17").Value = '30.706.44'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '''13.22'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '''0.000007531'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").Value = '2.141.33'
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("D22").Value = '''5.321'
$ws.Range("E22").Value = '  +1.59%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '''6.250'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("D25").Value = '''9.230'
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").Value = '''164.46'
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").Value = '''18.93'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '''1.922'
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("D29").Value = '''1.343'
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("D30").Value = '''0.09751'
$ws.Range("E30").Value = '  -2.84%  '
$ws.Range("D31").Value = '''1.493'
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").Value = '''4.303'
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("D33").Value = '''4.203'
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("D34").Value = '''0.04870'
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("D35").Value = '''1.128'
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D36").Value = '''0.6982'
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").Value = '''2.724'
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("D38").Value = '''0.01901'
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("D39").Value = '''2.806'
$ws.Range("E39").Value = '  +1.95%  '
$ws.Range("D40").Value = '''6.364'
$ws.Range("E40").Value = '  +1.07%  '
$ws.Range("D41").Value = '''75.98'
$ws.Range("E41").Value = '  +6.19%  '
$ws.Range("D42").Value = '''2.005'
$ws.Range("E42").Value = '  +1.62%  '
$ws.Range("D43").Value = '''0.4258'
$ws.Range("E43").Value = '  +1.83%  '
$ws.Range("D44").Value = '''1.000'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '''0.8364'
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '''101.68'
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("D47").Value = '''9.561'
$ws.Range("E47").Value = '  +2.80%  '
$ws.Range("D48").Value = '''35.71'
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("D49").Value = '''7.034'
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("D50").Value = '''918.91'
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").Value = '''0.05754'
